# ------------------------------------------------------------------
# "Implemented the title only"
#   * Slide 1 title/body get explicit Bold/Italic/Strikethrough/Underline
#     formatting, and the body text is replaced/extended.
#   * Slide 2 title gets explicit (disabled) Bold/Italic/Strikethrough/
#     Underline formatting.
#   * A new 3rd slide ("Title and Content" layout) is appended with a
#     title reusing the same look as Slide 1's title.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- Slide 1 : title "This is a title" -------------------------------
$s1 = $p.Slides.Item(1)

$titleRange1 = $s1.Shapes.Item(1).TextFrame.TextRange
$titleRange1.Font.Bold = $true
$titleRange1.Font.Italic = $true
$titleRange1.Font.Strikethrough = $false
$titleRange1.Font.Underline = $false

# --- Slide 1 : content placeholder "test" -----------------------------
$bodyRange1 = $s1.Shapes.Item(2).TextFrame.TextRange
$bodyRange1.Font.Bold = $false
$bodyRange1.Font.Italic = $false
$bodyRange1.Font.Strikethrough = $true
$bodyRange1.Font.Underline = $true
$bodyRange1.Text = "test adihwaudhwahid Decorations `n Yes sir"

# --- Slide 2 : title "Yes" --------------------------------------------
$s2 = $p.Slides.Item(2)

$titleRange2 = $s2.Shapes.Item(1).TextFrame.TextRange
$titleRange2.Font.Bold = $false
$titleRange2.Font.Italic = $false
$titleRange2.Font.Strikethrough = $false
$titleRange2.Font.Underline = $false

# --- Slide 3 (new) : "Title and Content" layout, title only ----------
$s3 = $p.Slides.Add(3, 2)

$titleRange3 = $s3.Shapes.Item(1).TextFrame.TextRange
$titleRange3.Text = "This is a List"
$titleRange3.Font.Name = "Calibri"
$titleRange3.Font.Size = 24
$titleRange3.Font.Color.RGB = 0
$titleRange3.Font.Bold = $true
$titleRange3.Font.Italic = $true
$titleRange3.Font.Strikethrough = $false
$titleRange3.Font.Underline = $false
